$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new "2022-Q1" sheet right before the "总计" sheet.
# ------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")
$new = $wb.Worksheets.Add($zj)
$new.Name = "2022-Q1"

# ------------------------------------------------------------------
# 2. Clone the cell formatting (border/bold/alignment) used by the
#    other quarterly sheets so the new sheet matches their look:
#    header row (B1:H1) + index column (A2:A15) use the bold/boxed
#    style, exactly like sheet "2021-Q4".
# ------------------------------------------------------------------
$src = $wb.Worksheets.Item("2021-Q4")
$src.Range("A1:H1").Copy()
$new.Range("A1:H1").PasteSpecial(-4122)
$src.Range("A2:H2").Copy()
$new.Range("A2:H15").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 3. Header row.
# ------------------------------------------------------------------
$new.Cells.Item(1,2).Value = "基金代码"
$new.Cells.Item(1,3).Value = "基金名称"
$new.Cells.Item(1,4).Value = "基金规模"
$new.Cells.Item(1,5).Value = "股票总仓位"
$new.Cells.Item(1,6).Value = "仓位占比"
$new.Cells.Item(1,7).Value = "持有市值(亿元)"
$new.Cells.Item(1,8).Value = "仓位排名"

# ------------------------------------------------------------------
# 4. Fund holdings data for 2022-Q1 (rows 2-15).
#    Numeric-looking text fields (fund code / scale / position /
#    ratio / market value) are entered with a leading "'" so they
#    stay text, matching the source data (t="inlineStr").
# ------------------------------------------------------------------
$new.Cells.Item(2,1).Value = 0
$new.Cells.Item(2,2).Value = "'166005"
$new.Cells.Item(2,3).Value = "中欧价值发现混合 -A"
$new.Cells.Item(2,4).Value = "'43.52"
$new.Cells.Item(2,5).Value = "'93.97"
$new.Cells.Item(2,6).Value = "'2.84"
$new.Cells.Item(2,7).Value = "'1.2360"
$new.Cells.Item(2,8).Value = 10

$new.Cells.Item(3,1).Value = 1
$new.Cells.Item(3,2).Value = "'001882"
$new.Cells.Item(3,3).Value = "中欧价值发现混合 -E"
$new.Cells.Item(3,4).Value = "'43.52"
$new.Cells.Item(3,5).Value = "'93.97"
$new.Cells.Item(3,6).Value = "'2.84"
$new.Cells.Item(3,7).Value = "'1.2360"
$new.Cells.Item(3,8).Value = 10

$new.Cells.Item(4,1).Value = 2
$new.Cells.Item(4,2).Value = "'001810"
$new.Cells.Item(4,3).Value = "中欧潜力价值灵活配置混合A"
$new.Cells.Item(4,4).Value = "'28.67"
$new.Cells.Item(4,5).Value = "'94.05"
$new.Cells.Item(4,6).Value = "'2.46"
$new.Cells.Item(4,7).Value = "'0.7053"
$new.Cells.Item(4,8).Value = 10

$new.Cells.Item(5,1).Value = 3
$new.Cells.Item(5,2).Value = "'090007"
$new.Cells.Item(5,3).Value = "大成策略回报混合"
$new.Cells.Item(5,4).Value = "'11.30"
$new.Cells.Item(5,5).Value = "'73.51"
$new.Cells.Item(5,6).Value = "'5.33"
$new.Cells.Item(5,7).Value = "'0.6023"
$new.Cells.Item(5,8).Value = 3

$new.Cells.Item(6,1).Value = 4
$new.Cells.Item(6,2).Value = "'008269"
$new.Cells.Item(6,3).Value = "大成睿享混合A"
$new.Cells.Item(6,4).Value = "'17.69"
$new.Cells.Item(6,5).Value = "'65.25"
$new.Cells.Item(6,6).Value = "'3.04"
$new.Cells.Item(6,7).Value = "'0.5378"
$new.Cells.Item(6,8).Value = 7

$new.Cells.Item(7,1).Value = 5
$new.Cells.Item(7,2).Value = "'004232"
$new.Cells.Item(7,3).Value = "中欧价值发现混合 -C"
$new.Cells.Item(7,4).Value = "'10.98"
$new.Cells.Item(7,5).Value = "'93.97"
$new.Cells.Item(7,6).Value = "'2.84"
$new.Cells.Item(7,7).Value = "'0.3118"
$new.Cells.Item(7,8).Value = 10

$new.Cells.Item(8,1).Value = 6
$new.Cells.Item(8,2).Value = "'011834"
$new.Cells.Item(8,3).Value = "大成投资严选六个月持有期混合型证券投资基金A"
$new.Cells.Item(8,4).Value = "'3.88"
$new.Cells.Item(8,5).Value = "'84.63"
$new.Cells.Item(8,6).Value = "'4.87"
$new.Cells.Item(8,7).Value = "'0.1890"
$new.Cells.Item(8,8).Value = 8

$new.Cells.Item(9,1).Value = 7
$new.Cells.Item(9,2).Value = "'166024"
$new.Cells.Item(9,3).Value = "中欧恒利三年定期开放混合"
$new.Cells.Item(9,4).Value = "'4.48"
$new.Cells.Item(9,5).Value = "'98.71"
$new.Cells.Item(9,6).Value = "'3.69"
$new.Cells.Item(9,7).Value = "'0.1653"
$new.Cells.Item(9,8).Value = 8

$new.Cells.Item(10,1).Value = 8
$new.Cells.Item(10,2).Value = "'013463"
$new.Cells.Item(10,3).Value = "大成致远优势一年持有期混合A"
$new.Cells.Item(10,4).Value = "'4.01"
$new.Cells.Item(10,5).Value = "'60.15"
$new.Cells.Item(10,6).Value = "'4.04"
$new.Cells.Item(10,7).Value = "'0.1620"
$new.Cells.Item(10,8).Value = 7

$new.Cells.Item(11,1).Value = 9
$new.Cells.Item(11,2).Value = "'090013"
$new.Cells.Item(11,3).Value = "大成竞争优势混合"
$new.Cells.Item(11,4).Value = "'3.87"
$new.Cells.Item(11,5).Value = "'72.16"
$new.Cells.Item(11,6).Value = "'3.80"
$new.Cells.Item(11,7).Value = "'0.1471"
$new.Cells.Item(11,8).Value = 7

$new.Cells.Item(12,1).Value = 10
$new.Cells.Item(12,2).Value = "'008270"
$new.Cells.Item(12,3).Value = "大成睿享混合C"
$new.Cells.Item(12,4).Value = "'2.87"
$new.Cells.Item(12,5).Value = "'65.25"
$new.Cells.Item(12,6).Value = "'3.04"
$new.Cells.Item(12,7).Value = "'0.0872"
$new.Cells.Item(12,8).Value = 7

$new.Cells.Item(13,1).Value = 11
$new.Cells.Item(13,2).Value = "'005764"
$new.Cells.Item(13,3).Value = "中欧潜力价值灵活配置混合C"
$new.Cells.Item(13,4).Value = "'3.43"
$new.Cells.Item(13,5).Value = "'94.05"
$new.Cells.Item(13,6).Value = "'2.46"
$new.Cells.Item(13,7).Value = "'0.0844"
$new.Cells.Item(13,8).Value = 10

$new.Cells.Item(14,1).Value = 12
$new.Cells.Item(14,2).Value = "'011835"
$new.Cells.Item(14,3).Value = "大成投资严选六个月持有期混合型证券投资基金C"
$new.Cells.Item(14,4).Value = "'0.30"
$new.Cells.Item(14,5).Value = "'84.63"
$new.Cells.Item(14,6).Value = "'4.87"
$new.Cells.Item(14,7).Value = "'0.0146"
$new.Cells.Item(14,8).Value = 8

$new.Cells.Item(15,1).Value = 13
$new.Cells.Item(15,2).Value = "'013464"
$new.Cells.Item(15,3).Value = "大成致远优势一年持有期混合C"
$new.Cells.Item(15,4).Value = "'0.17"
$new.Cells.Item(15,5).Value = "'60.15"
$new.Cells.Item(15,6).Value = "'4.04"
$new.Cells.Item(15,7).Value = "'0.0069"
$new.Cells.Item(15,8).Value = 7

# ------------------------------------------------------------------
# 5. Update the "总计" (summary) sheet: insert a new 2022-Q1 row
#    right after the header, pushing the existing quarters down by
#    one row and bumping their index-column (A) counters by 1.
# ------------------------------------------------------------------
$oldB = @("", "", "", "")
$oldC = @(0, 0, 0, 0)
$oldD = @(0, 0, 0, 0)
for ($i = 0; $i -lt 4; $i++) {
    $r = 2 + $i
    $oldB[$i] = $zj.Cells.Item($r, 2).Value()
    $oldC[$i] = $zj.Cells.Item($r, 3).Value()
    $oldD[$i] = $zj.Cells.Item($r, 4).Value()
}

# Extend the boxed index-column style down to the new last row (A6).
$zj.Range("A5").Copy()
$zj.Range("A6").PasteSpecial(-4122)

for ($i = 3; $i -ge 0; $i--) {
    $r = 3 + $i
    $zj.Cells.Item($r, 1).Value = $i + 1
    $zj.Cells.Item($r, 2).Value = $oldB[$i]
    $zj.Cells.Item($r, 3).Value = $oldC[$i]
    $zj.Cells.Item($r, 4).Value = $oldD[$i]
}

$zj.Cells.Item(2, 1).Value = 0
$zj.Cells.Item(2, 2).Value = "2022-Q1"
$zj.Cells.Item(2, 3).Value = 14
$zj.Cells.Item(2, 4).Value = 5.49
